# Apply corrected error-estimation / projected-years results to the
# "Trends Status" and "Species qualification" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$wsTrends = $wb.Worksheets.Item("Trends Status")

$wsTrends.Range("C2").Value = 3
$wsTrends.Range("D2").Value = 33.3
$wsTrends.Range("E2").Value = 10.7

$wsTrends.Range("B3").Value = 6
$wsTrends.Range("C3").Value = 10
$wsTrends.Range("D3").Value = 40
$wsTrends.Range("E3").Value = 35.7

$wsTrends.Range("B4").Value = 3
$wsTrends.Range("C4").Value = 11
$wsTrends.Range("D4").Value = 20
$wsTrends.Range("E4").Value = 39.3

$wsTrends.Range("C5").Value = 3
$wsTrends.Range("E5").Value = 10.7

$wsTrends.Range("D6").Value = 6.7
$wsTrends.Range("E6").Value = 3.6

$wsTrends.Range("B7").Value = 16
$wsTrends.Range("C7").Value = 29

# --- Sheet: "Species qualification" ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")

$wsSpecies.Range("C3").Value = 15
$wsSpecies.Range("C4").Value = 28
